# Fish Bulletin 163 - Table22: revisions on partyboat landings
# Fill in the previously-blank "Number of anglers" row (row 22) with the
# revised angler counts for each year 1964-1973 (columns B-K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B22").Value = 695445
$ws.Range("C22").Value = 688103
$ws.Range("D22").Value = 857000
$ws.Range("E22").Value = 780102
$ws.Range("F22").Value = 849654
$ws.Range("G22").Value = 802811
$ws.Range("H22").Value = 872327
$ws.Range("I22").Value = 728126
$ws.Range("J22").Value = 792618
$ws.Range("K22").Value = 880100

# Reflect the updated view state captured in the saved workbook: the author
# scrolled right/zoomed out a bit and left the cursor on I10 after editing.
$win = $excel.ActiveWindow
$win.Zoom = 120
$ws.Range("I10").Select() | Out-Null
